$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values that apply identically to rows 2 and 3
$values = @{
    "I" = 0
    "J" = 0
    "K" = 46.9
    "L" = 0.2720417633410673
    "M" = 21.6
    "N" = 0.1520056298381422
    "O" = 0.4605543710021323
    "P" = 21.6
    "Q" = 0.1520056298381422
    "R" = 0.4605543710021323
    "U" = 52.9
    "V" = 0.3722730471498945
    "W" = 0.1242713301536831
    "X" = 0.03901251200783102
    "Y" = 0.08525881814585207
    "Z" = 0.2199007640403577
    "AA" = 0
    "AB" = 0.03173547195027669
    "AC" = -0.03173547195027669
    "AD" = 412.7
    "AE" = 0
    "AF" = 412.7
    "AG" = 359.8
    "AH" = 0.7438716654650325
    "AI" = 0.5183371012308465
    "AJ" = 0.7168758716875873
    "AK" = 0.4840575810574466
}

foreach ($row in 2,3) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
    # Remove debt_ebitda (AN) and net_debt_ebitda (AP) values entirely
    $ws.Range("AN$row").ClearContents()
    $ws.Range("AP$row").ClearContents()
}
